# Update report data for Montefiorino through 2021-12-08 ("aggiornamento fino a 8/12").
# Extends the daily series in Sheet1 from row 385 (2021-09-20) through row 464 (2021-12-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new rows (386:464) with the same style/formatting as the last existing
# data row (385) so column A keeps its date-style (numFmtId 165) formatting.
$srcRange = $ws.Range("A385:D385")
$dstRange = $ws.Range("A386:D464")
$srcRange.Copy($dstRange)

# row, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
  @(386, 44460, 0, 0, 0),
  @(387, 44461, 0, 0, 0),
  @(388, 44462, 0, 0, 0),
  @(389, 44463, 0, 0, 0),
  @(390, 44464, 0, 0, 0),
  @(391, 44465, 0, 0, 0),
  @(392, 44466, 0, 0, 0),
  @(393, 44467, 0, 0, 0),
  @(394, 44468, 0, 0, 0),
  @(395, 44469, 0, 0, 0),
  @(396, 44470, 0, 0, 0),
  @(397, 44471, 0, 0, 0),
  @(398, 44472, 0, 0, 0),
  @(399, 44473, 0, 0, 0),
  @(400, 44474, 0, 0, 0),
  @(401, 44475, 0, 0, 0),
  @(402, 44476, 0, 0, 0),
  @(403, 44477, 0, 0, 0),
  @(404, 44478, 0, 0, 0),
  @(405, 44479, 0, 0, 0),
  @(406, 44480, 0, 0, 0),
  @(407, 44481, 0, 0, 0),
  @(408, 44482, 0, 0, 0),
  @(409, 44483, 0, 0, 0),
  @(410, 44484, 0, 0, 0),
  @(411, 44485, 0, 0, 0),
  @(412, 44486, 0, 0, 0),
  @(413, 44487, 0, 0, 0),
  @(414, 44488, 0, 0, 0),
  @(415, 44489, 0, 0, 0),
  @(416, 44490, 0, 0, 0),
  @(417, 44491, 0, 0, 0),
  @(418, 44492, 0, 0, 0),
  @(419, 44493, 0, 0, 0),
  @(420, 44494, 0, 0, 0),
  @(421, 44495, 0, 0, 0),
  @(422, 44496, 0, 0, 0),
  @(423, 44497, 0, 0, 0),
  @(424, 44498, 0, 0, 0),
  @(425, 44499, 0, 0, 0),
  @(426, 44500, 0, 0, 0),
  @(427, 44501, 0, 0, 0),
  @(428, 44502, 0, 0, 0),
  @(429, 44503, 0, 0, 0),
  @(430, 44504, 0, 0, 0),
  @(431, 44505, 0, 0, 0),
  @(432, 44506, 0, 0, 0),
  @(433, 44507, 0, 0, 0),
  @(434, 44508, 0, 0, 0),
  @(435, 44509, 0, 0, 0),
  @(436, 44510, 0, 0, 0),
  @(437, 44511, 0, 0, 0),
  @(438, 44512, 0, 0, 0),
  @(439, 44513, 0, 0, 0),
  @(440, 44514, 0, 0, 0),
  @(441, 44515, 0, 0, 0),
  @(442, 44516, 1, 1, 46.70714619336758),
  @(443, 44517, 0, 1, 46.70714619336758),
  @(444, 44518, 0, 1, 46.70714619336758),
  @(445, 44519, 0, 1, 46.70714619336758),
  @(446, 44520, 0, 1, 46.70714619336758),
  @(447, 44521, 0, 1, 46.70714619336758),
  @(448, 44522, 0, 1, 46.70714619336758),
  @(449, 44523, 0, 0, 0),
  @(450, 44524, 1, 1, 46.70714619336758),
  @(451, 44525, 0, 1, 46.70714619336758),
  @(452, 44526, 0, 1, 46.70714619336758),
  @(453, 44527, 0, 1, 46.70714619336758),
  @(454, 44528, 0, 1, 46.70714619336758),
  @(455, 44529, 0, 1, 46.70714619336758),
  @(456, 44530, 0, 1, 46.70714619336758),
  @(457, 44531, 0, 0, 0),
  @(458, 44532, 0, 0, 0),
  @(459, 44533, 1, 1, 46.70714619336758),
  @(460, 44534, 0, 1, 46.70714619336758),
  @(461, 44535, 0, 1, 46.70714619336758),
  @(462, 44536, 0, 1, 46.70714619336758),
  @(463, 44537, 0, 1, 46.70714619336758),
  @(464, 44538, 0, 1, 46.70714619336758)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Output "Added rows 386-464 (through 2021-12-08) to Sheet1."
